$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) date value from 45204 to 45205 for all data rows (2..530)
$ws.Range("C2:C530").Value = 45205
